$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 306, pushing the
# existing rows 306-330 down to 308-332 (preserving their formatting,
# including the date style on column D).
$ws.Rows("306:307").Insert()

# --- New row 306: Femacal de La Calera / Repollo / Crespo record / Primera ---
$ws.Cells.Item(306, 1).Value = 3
$ws.Cells.Item(306, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(306, 3).Value = "Coquimbo"
$ws.Cells.Item(306, 4).Value = 44461
$ws.Cells.Item(306, 4).NumberFormat = $ws.Cells.Item(305, 4).NumberFormat()
$ws.Cells.Item(306, 5).Value = 5
$ws.Cells.Item(306, 6).Value = 100112006
$ws.Cells.Item(306, 7).Value = "Repollo"
$ws.Cells.Item(306, 8).Value = "Crespo record"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 1200
$ws.Cells.Item(306, 11).Value = 600
$ws.Cells.Item(306, 12).Value = 600
$ws.Cells.Item(306, 13).Value = 600
$ws.Cells.Item(306, 14).Value = "`$/unidad"
$ws.Cells.Item(306, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(306, 16).Value = 600
$ws.Cells.Item(306, 17).Value = 1
$ws.Cells.Item(306, 18).Value = "Hortaliza"

# --- New row 307: Femacal de La Calera / Repollo / Crespo record / Segunda ---
$ws.Cells.Item(307, 1).Value = 3
$ws.Cells.Item(307, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(307, 3).Value = "Coquimbo"
$ws.Cells.Item(307, 4).Value = 44461
$ws.Cells.Item(307, 4).NumberFormat = $ws.Cells.Item(305, 4).NumberFormat()
$ws.Cells.Item(307, 5).Value = 5
$ws.Cells.Item(307, 6).Value = 100112006
$ws.Cells.Item(307, 7).Value = "Repollo"
$ws.Cells.Item(307, 8).Value = "Crespo record"
$ws.Cells.Item(307, 9).Value = "Segunda"
$ws.Cells.Item(307, 10).Value = 880
$ws.Cells.Item(307, 11).Value = 500
$ws.Cells.Item(307, 12).Value = 500
$ws.Cells.Item(307, 13).Value = 500
$ws.Cells.Item(307, 14).Value = "`$/unidad"
$ws.Cells.Item(307, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(307, 16).Value = 500
$ws.Cells.Item(307, 17).Value = 1
$ws.Cells.Item(307, 18).Value = "Hortaliza"
